$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-4 and add new row 5 with corrected packet info
$ws.Range("A2").Value = "162.159.135.234"
$ws.Range("B2").Value = 443
$ws.Range("C2").Value = "Ether / IP / TCP 192.168.1.80:52490 > 162.159.135.234:https A"

$ws.Range("A3").Value = "185.199.111.154"
$ws.Range("B3").Value = 443
$ws.Range("C3").Value = "Ether / IP / TCP 192.168.1.80:52821 > 185.199.111.154:https A / Raw"

$ws.Range("A4").Value = "185.199.111.154"
$ws.Range("B4").Value = 443
$ws.Range("C4").Value = "Ether / IP / TCP 192.168.1.80:52813 > 185.199.111.154:https A / Raw"

$ws.Range("A5").Value = "140.82.114.26"
$ws.Range("B5").Value = 443
$ws.Range("C5").Value = "Ether / IP / TCP 192.168.1.80:52824 > 140.82.114.26:https A / Raw"
